# De-dupe the stray " - as percentages" rows from the 2002 Intermediate_2
# export (once in the Male block at row 68, once in the Female block at
# row 135). Deleting these two rows shifts everything below each of them
# up by one, matching the "fixed 2002, de dupe is good for test output
# (csv)" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the later row first so the earlier row's index ("68") stays valid.
$ws.Rows("135").Delete() | Out-Null
$ws.Rows("68").Delete() | Out-Null
